$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 587.5
$ws.Range("I2").Value = 587.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 587.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -474.5
$ws.Range("N2").Value = $null

$ws.Range("H9").Value = 1278.6
$ws.Range("I9").Value = 550
$ws.Range("J9").Value = 1764.3334
$ws.Range("K9").Value = 550
$ws.Range("L9").Value = 1764.3334
$ws.Range("M9").Value = -381
$ws.Range("N9").Value = -2102.3334

$ws.Range("H18").Value = 1065.75
$ws.Range("I18").Value = 980.9091
$ws.Range("K18").Value = 980.9091
$ws.Range("M18").Value = -696.9091

$ws.Range("H19").Value = 599.75
$ws.Range("I19").Value = 599.75
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 599.75
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -424.75
$ws.Range("N19").Value = $null

$ws.Range("H28").Value = 447.45456
$ws.Range("I28").Value = 102.77778
$ws.Range("K28").Value = 102.77778
$ws.Range("M28").Value = 382.22222

$ws.Range("H38").Value = 163.71428
$ws.Range("I38").Value = 24.5
$ws.Range("K38").Value = 73.5
$ws.Range("M38").Value = 298.5

$ws.Range("H106").Value = 200002610
$ws.Range("I106").Value = 200002610
$ws.Range("K106").Value = 200002610
$ws.Range("M106").Value = -200001979

$ws.Range("H107").Value = 64761.855
$ws.Range("I107").Value = 82314.63
$ws.Range("J107").Value = 401.66666
$ws.Range("K107").Value = 82314.63
$ws.Range("L107").Value = 401.66666
$ws.Range("M107").Value = -80394.63
$ws.Range("N107").Value = -4241.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 93
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 93
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 93
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = -319

$ws.Range("H45").Value = 2500
$ws.Range("I45").Value = 2500
$ws.Range("K45").Value = 2500
$ws.Range("M45").Value = -2123

$ws.Range("H110").Value = 50000900
$ws.Range("I110").Value = 502.66666
$ws.Range("J110").Value = 125001500
$ws.Range("K110").Value = 502.66666
$ws.Range("L110").Value = 125001500
$ws.Range("M110").Value = 1542.33334
$ws.Range("N110").Value = -125005590

$ws.Range("H116").Value = 93
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 93
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 93
$ws.Range("M116").Value = $null
$ws.Range("N116").Value = -4681

$ws.Range("H122").Value = 1337
$ws.Range("I122").Value = 1337
$ws.Range("K122").Value = 4011
$ws.Range("M122").Value = -1561

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 93
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 93
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 93
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -321

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null

$ws.Range("H99").Value = 2326.4736
$ws.Range("I99").Value = 2326.4736
$ws.Range("K99").Value = 2326.4736
$ws.Range("M99").Value = -828.4735999999998

$ws.Range("H105").Value = 2010
$ws.Range("J105").Value = 2010
$ws.Range("L105").Value = 2010
$ws.Range("N105").Value = -5504

$ws.Range("H107").Value = 135663.33
$ws.Range("I107").Value = 201745
$ws.Range("K107").Value = 201745
$ws.Range("M107").Value = -199825

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 250050
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 12

$ws.Range("H16").Value = 1395
$ws.Range("I16").Value = 1245.625
$ws.Range("J16").Value = 1992.5
$ws.Range("K16").Value = 1245.625
$ws.Range("L16").Value = 1992.5
$ws.Range("M16").Value = -958.625
$ws.Range("N16").Value = -2566.5

$ws.Range("H105").Value = 626.4286
$ws.Range("I105").Value = 564.1667
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 564.1667
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 1182.8333
$ws.Range("N105").Value = -4494

$ws.Range("H107").Value = 750
$ws.Range("I107").Value = 750
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 750
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1170
$ws.Range("N107").Value = $null

$ws.Range("H113").Value = 1395
$ws.Range("I113").Value = 1245.625
$ws.Range("J113").Value = 1992.5
$ws.Range("K113").Value = 1245.625
$ws.Range("L113").Value = 1992.5
$ws.Range("M113").Value = 924.375
$ws.Range("N113").Value = -6332.5

$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 150000
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 52497.25
$ws.Range("I139").Value = 4994.5
$ws.Range("K139").Value = 14983.5
$ws.Range("M139").Value = -9843.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null

$ws.Range("H80").Value = 5612
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5612
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 5612
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = -7608

$ws.Range("H83").Value = 5612
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5612
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 28060
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = -38044

$ws.Range("H97").Value = 281
$ws.Range("I97").Value = 296.5
$ws.Range("J97").Value = 250
$ws.Range("K97").Value = 296.5
$ws.Range("L97").Value = 250
$ws.Range("M97").Value = 199.5
$ws.Range("N97").Value = -1242

$ws.Range("H113").Value = 750
$ws.Range("I113").Value = 750
$ws.Range("K113").Value = 750
$ws.Range("M113").Value = 1420

$ws.Range("H122").Value = 5103.154
$ws.Range("I122").Value = 3855.75
$ws.Range("K122").Value = 11567.25
$ws.Range("M122").Value = -9117.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 25
$ws.Range("I2").Value = 25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 87
$ws.Range("N2").Value = $null

$ws.Range("H17").Value = 256
$ws.Range("I17").Value = 134
$ws.Range("K17").Value = 134
$ws.Range("M17").Value = 36

$ws.Range("H19").Value = 361
$ws.Range("I19").Value = 361
$ws.Range("K19").Value = 361
$ws.Range("M19").Value = -191

$ws.Range("H22").Value = 1355.85
$ws.Range("I22").Value = 1044.1
$ws.Range("J22").Value = 1667.6
$ws.Range("K22").Value = 1044.1
$ws.Range("L22").Value = 1667.6
$ws.Range("M22").Value = -749.0999999999999
$ws.Range("N22").Value = -2257.6

$ws.Range("H23").Value = 3000
$ws.Range("I23").Value = 3000
$ws.Range("K23").Value = 3000
$ws.Range("M23").Value = -2770

$ws.Range("H27").Value = 1355.85
$ws.Range("I27").Value = 1044.1
$ws.Range("J27").Value = 1667.6
$ws.Range("K27").Value = 1044.1
$ws.Range("L27").Value = 1667.6
$ws.Range("M27").Value = -937.0999999999999
$ws.Range("N27").Value = -1881.6

$ws.Range("H61").Value = 4752
$ws.Range("I61").Value = 4002.6667
$ws.Range("K61").Value = 4002.6667
$ws.Range("M61").Value = -3800.6667

$ws.Range("H82").Value = 3205.4614
$ws.Range("I82").Value = 2097.4285
$ws.Range("J82").Value = 4498.1665
$ws.Range("K82").Value = 2097.4285
$ws.Range("L82").Value = 4498.1665
$ws.Range("M82").Value = -1736.4285
$ws.Range("N82").Value = -5220.1665

$ws.Range("H85").Value = 3205.4614
$ws.Range("I85").Value = 2097.4285
$ws.Range("J85").Value = 4498.1665
$ws.Range("K85").Value = 2097.4285
$ws.Range("L85").Value = 4498.1665
$ws.Range("M85").Value = -849.4285
$ws.Range("N85").Value = -6994.1665

$ws.Range("H113").Value = 4752
$ws.Range("I113").Value = 4002.6667
$ws.Range("K113").Value = 4002.6667
$ws.Range("M113").Value = -1832.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 350.125
$ws.Range("I100").Value = 257.2857
$ws.Range("K100").Value = 514.5714
$ws.Range("M100").Value = 26.42859999999996

$ws.Range("H113").Value = 242.625
$ws.Range("I113").Value = 191
$ws.Range("J113").Value = 397.5
$ws.Range("K113").Value = 573
$ws.Range("L113").Value = 1192.5
$ws.Range("M113").Value = 1597
$ws.Range("N113").Value = -5532.5

$ws.Range("H122").Value = 1256
$ws.Range("I122").Value = 958.6
$ws.Range("K122").Value = 2875.8
$ws.Range("M122").Value = -425.8000000000002

